# Input_Filelist.xlsx update: add newly-computed rows (But-OH, Benz-OH,
# TMM-OH, TMA-OH, W-OH, MIm-OH) and extra temperature columns (475/525 K)
# to the fitting results table, per "Add files via upload".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 1 (temperature header row): reshuffle columns so a 4.75 / 5.25
# column is inserted and the whole row shifts right by two columns.
# F1/G1 need the "boxed" style that currently lives on I1 (s=2); the
# rest of the row keeps the plain bordered style that currently lives
# on E1/H1 (s=1). Copy the formatting *before* overwriting the values
# so the source cells are still intact when we copy from them.
# ------------------------------------------------------------------
$ws.Range("I1").Copy() | Out-Null
$ws.Range("F1:G1").PasteSpecial(-4122) | Out-Null

$ws.Range("E1").Copy() | Out-Null
$ws.Range("I1:K1").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

$ws.Range("E1").Value = 4.75
$ws.Range("F1").Value = 5
$ws.Range("G1").Value = 5.25
$ws.Range("H1").Value = 5.5
$ws.Range("I1").Value = 6
$ws.Range("J1").Value = 6.5
$ws.Range("K1").Value = 7

# ------------------------------------------------------------------
# Row 3's H3 used to carry the wrap-text "Explanatory Text" style
# (with the matching ht="72" row override). The refreshed table no
# longer wraps any text, so clear that formatting back to plain
# before the new values go in (copy the neighbouring plain cell's
# format so no brand-new style slot gets allocated).
# ------------------------------------------------------------------
$ws.Range("G3").Copy() | Out-Null
$ws.Range("H3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Rows.Item(3).AutoFit() | Out-Null

# ------------------------------------------------------------------
# Data rows 2-8: one row per ligand (OH-OH, But-OH, Benz-OH, TMM-OH,
# TMA-OH, W-OH, MIm-OH), columns A-D unchanged layout, columns E-K
# now hold the 475/500/525/550/600/650/700 K log-file paths.
# ------------------------------------------------------------------

# Row 2 - OH-OH
$ws.Range("A2").Value = "OH-OH"
$ws.Range("B2").Value = -305.10889826
$ws.Range("C2").Value = -305.10889826
$ws.Range("D2").Value = 200
$ws.Range("E2").Value = "../New_OH-OH/OH-OH_475.log"
$ws.Range("F2").Value = "../New_OH-OH/OH-OH_500.log"
$ws.Range("G2").Value = "../New_OH-OH/OH-OH_525.log"
$ws.Range("H2").Value = "../New_OH-OH/OH-OH_550.log"
$ws.Range("I2").Value = "../New_OH-OH/OH-OH_600.log"
$ws.Range("J2").Value = "../New_OH-OH/OH-OH_650.log"
$ws.Range("K2").Value = "../New_OH-OH/OH-OH_700.log"

# Row 3 - But-OH
$ws.Range("A3").Value = "But-OH"
$ws.Range("B3").Value = -158.45799844
$ws.Range("C3").Value = -305.10889826
$ws.Range("D3").Value = 200
$ws.Range("E3").Value = "../New_But-OH/But_OH_475.log"
$ws.Range("F3").Value = "../New_But-OH/But_OH_500.log"
$ws.Range("G3").Value = "../New_But-OH/But_OH_525.log"
$ws.Range("H3").Value = "../New_But-OH/But_OH_550.log"
$ws.Range("I3").Value = "../New_But-OH/But_OH_600.log"
$ws.Range("J3").Value = "../New_But-OH/But_OH_650.log"
$ws.Range("K3").Value = "../New_But-OH/But_OH_700.log"

# Row 4 - Benz-OH
$ws.Range("A4").Value = "Benz-OH"
$ws.Range("B4").Value = -232.24858806
$ws.Range("C4").Value = -305.10889826
$ws.Range("D4").Value = 180
$ws.Range("E4").Value = "../New_Benz-OH/Benz_OH_475.log"
$ws.Range("F4").Value = "../New_Benz-OH/Benz_OH_500.log"
$ws.Range("G4").Value = "../New_Benz-OH/Benz_OH_525.log"
$ws.Range("H4").Value = "../New_Benz-OH/Benz_OH_550.log"
$ws.Range("I4").Value = "../New_Benz-OH/Benz_OH_600.log"
$ws.Range("J4").Value = "../New_Benz-OH/Benz_OH_650.log"
$ws.Range("K4").Value = "../New_Benz-OH/Benz_OH_700.log"

# Row 5 - TMM-OH
$ws.Range("A5").Value = "TMM-OH"
$ws.Range("B5").Value = -158.45877153
$ws.Range("C5").Value = -305.10889826
$ws.Range("D5").Value = 200
$ws.Range("E5").Value = "../New_TMM-OH/TMM_OH_475.log"
$ws.Range("F5").Value = "../New_TMM-OH/TMM_OH_500.log"
$ws.Range("G5").Value = "../New_TMM-OH/TMM_OH_525.log"
$ws.Range("H5").Value = "../New_TMM-OH/TMM_OH_550.log"
$ws.Range("I5").Value = "../New_TMM-OH/TMM_OH_600.log"
$ws.Range("J5").Value = "../New_TMM-OH/TMM_OH_650.log"
$ws.Range("K5").Value = "../New_TMM-OH/TMM_OH_700.log"

# Row 6 - TMA-OH
$ws.Range("A6").Value = "TMA-OH"
$ws.Range("B6").Value = -174.85101956
$ws.Range("C6").Value = -305.10889826
$ws.Range("D6").Value = 110
$ws.Range("E6").Value = "../New_TMA-OH/TMA_OH_475.log"
$ws.Range("F6").Value = "../New_TMA-OH/TMA_OH_500.log"
$ws.Range("G6").Value = "../New_TMA-OH/TMA_OH_525.log"
$ws.Range("H6").Value = "../New_TMA-OH/TMA_OH_550.log"
$ws.Range("I6").Value = "../New_TMA-OH/TMA_OH_600.log"
$ws.Range("J6").Value = "../New_TMA-OH/TMA_OH_650.log"
$ws.Range("K6").Value = "../New_TMA-OH/TMA_OH_700.log"

# Row 7 - W-OH
$ws.Range("A7").Value = "W-OH"
$ws.Range("B7").Value = -305.68926621
$ws.Range("C7").Value = -305.10889826
$ws.Range("D7").Value = 200
$ws.Range("E7").Value = "../New_W-OH/W_OH_475.log"
$ws.Range("F7").Value = "../New_W-OH/W_OH_500.log"
$ws.Range("G7").Value = "../New_W-OH/W_OH_525.log"
$ws.Range("H7").Value = "../New_W-OH/W_OH_550.log"
$ws.Range("I7").Value = "../New_W-OH/W_OH_600.log"
$ws.Range("J7").Value = "../New_W-OH/W_OH_650.log"
$ws.Range("K7").Value = "../New_W-OH/W_OH_700.log"

# Row 8 - MIm-OH (no 475 K log for this ligand, so E8 stays empty)
$ws.Range("A8").Value = "MIm-OH"
$ws.Range("B8").Value = -265.923572127
$ws.Range("C8").Value = -305.10889826
$ws.Range("D8").Value = 110
$ws.Range("F8").Value = "../New_MIm-OH/MIm_OH_500.log"
$ws.Range("G8").Value = "../New_MIm-OH/MIm_OH_525.log"
$ws.Range("H8").Value = "../New_MIm-OH/MIm_OH_550.log"
$ws.Range("I8").Value = "../New_MIm-OH/MIm_OH_600.log"
$ws.Range("J8").Value = "../New_MIm-OH/MIm_OH_650.log"
$ws.Range("K8").Value = "../New_MIm-OH/MIm_OH_700.log"

# ------------------------------------------------------------------
# Widen the new log-path / label columns so the longer strings are
# readable (characters units, same as Format > Column Width).
# ------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 12.6640625
$ws.Columns.Item(6).ColumnWidth = 25.5546875
$ws.Columns.Item(7).ColumnWidth = 17.44140625
$ws.Columns.Item(8).ColumnWidth = 17.21875
$ws.Columns.Item(9).ColumnWidth = 17.88671875
$ws.Columns.Item(10).ColumnWidth = 21.88671875
$ws.Columns.Item(11).ColumnWidth = 31.6640625

# ------------------------------------------------------------------
# Leave the selection where the author last left it.
# ------------------------------------------------------------------
$ws.Range("F10").Select() | Out-Null
